$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new rows at row 410; existing rows 410-455 shift down to 414-459 ---
$ws.Rows("410:413").Insert()

# --- A:C and E:J are identical "catalog" columns on every data row in this sheet ---
$abc = New-Object 'object[,]' 4,3
$ej  = New-Object 'object[,]' 4,6
for ($i = 0; $i -lt 4; $i++) {
    $abc[$i,0] = 5
    $abc[$i,1] = "Macroferia Regional de Talca"
    $abc[$i,2] = "Maule"
    $ej[$i,0] = 7
    $ej[$i,1] = "Fruta"
    $ej[$i,2] = 100103
    $ej[$i,3] = "Frutos de hueso (carozo)"
    $ej[$i,4] = 100103004
    $ej[$i,5] = "Durazno"
}
$ws.Range("A410:C413").Value = $abc
$ws.Range("E410:J413").Value = $ej

# --- Column D (Fecha, date-formatted) ---
$dcol = New-Object 'object[,]' 4,1
$dcol[0,0] = 44918
$dcol[1,0] = 44918
$dcol[2,0] = 44918
$dcol[3,0] = 44918
$ws.Range("D410:D413").Value = $dcol

# --- Columns K:T (Variedad..Kg/unidad) for the 4 new rows ---
$kt = New-Object 'object[,]' 4,10
$kt[0,0] = "Kurakata"
$kt[0,1] = "Extra (doble especial)"
$kt[0,2] = 150
$kt[0,3] = 18000
$kt[0,4] = 18000
$kt[0,5] = 18000
$kt[0,6] = "`$/bandeja 15 kilos granel"
$kt[0,7] = "Región de O'Higgins"
$kt[0,8] = 1200
$kt[0,9] = 15

$kt[1,0] = "Royal Glory"
$kt[1,1] = "Especial"
$kt[1,2] = 250
$kt[1,3] = 12000
$kt[1,4] = 12000
$kt[1,5] = 12000
$kt[1,6] = "`$/bandeja 15 kilos granel"
$kt[1,7] = "Región de O'Higgins"
$kt[1,8] = 800
$kt[1,9] = 15

$kt[2,0] = "Royal Glory"
$kt[2,1] = "Primera"
$kt[2,2] = 200
$kt[2,3] = 10000
$kt[2,4] = 10000
$kt[2,5] = 10000
$kt[2,6] = "`$/bandeja 15 kilos granel"
$kt[2,7] = "Región de O'Higgins"
$kt[2,8] = 667
$kt[2,9] = 15

$kt[3,0] = "Royal Glory"
$kt[3,1] = "Especial"
$kt[3,2] = 180
$kt[3,3] = 8000
$kt[3,4] = 8000
$kt[3,5] = 8000
$kt[3,6] = "`$/bandeja 15 kilos granel"
$kt[3,7] = "Región de O'Higgins"
$kt[3,8] = 533
$kt[3,9] = 15

$ws.Range("K410:T413").Value = $kt
